$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for 2018-05-23, appended after existing row 18 (rows 19-24)
$data = @(
    @("2018-05-23", "161725", "招商中证白酒指数分级", "1.2100", "1.1945", "-1.28"),
    @("2018-05-23", "000311", "景顺长城沪深300增强", "2.0560", "2.0396", "-0.80"),
    @("2018-05-23", "110022", "易方达消费行业", "2.3570", "2.3327", "-1.03"),
    @("2018-05-23", "161616", "融通医疗保健行业混合", "1.2860", "1.2878", "0.14"),
    @("2018-05-23", "486001", "工银全球配置", "1.3020", "1.2971", "-0.37"),
    @("2018-05-23", "002611", "博时黄金ETF联接C", "0.9349", "0.9363", "0.15")
)

$startRow = 19
$endRow = $startRow + $data.Count - 1
$fullRange = $ws.Range("A" + $startRow + ":F" + $endRow)

# Force text formatting first so values like dates, fund codes and
# decimal-looking strings are stored verbatim instead of being
# auto-converted to numbers/dates.
$fullRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
}

# Drop the explicit text style again so the cells fall back to the
# workbook's default (unstyled) formatting, matching the rest of the sheet.
$fullRange.ClearFormats()
